# Updates to team matrix values (Canisius_B) from games pulled March 7.
# Only the probability cells that changed in the source commit are touched;
# every other cell in the sheet is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1985294117647059
$ws.Range("C2").Value = 0.5330882352941176
$ws.Range("J2").Value = 0.01470588235294118
$ws.Range("P2").Value = 0.1544117647058824
$ws.Range("S2").Value = 0.09926470588235294

# Row 3
$ws.Range("B3").Value = 0.006666666666666667
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.02
$ws.Range("P3").Value = 0.72
$ws.Range("S3").Value = 0.22

# Row 6
$ws.Range("B6").Value = 0.04700854700854701
$ws.Range("D6").Value = 0.008547008547008548
$ws.Range("F6").Value = 0.05982905982905983
$ws.Range("J6").Value = 0.2863247863247863
$ws.Range("O6").Value = 0.02991452991452992
$ws.Range("Q6").Value = 0.1837606837606838
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.3290598290598291

# Row 7
$ws.Range("B7").Value = 0.09433962264150944
$ws.Range("D7").Value = 0.03773584905660377
$ws.Range("F7").Value = 0.06918238993710692
$ws.Range("J7").Value = 0.2012578616352201
$ws.Range("O7").Value = 0.01257861635220126
$ws.Range("Q7").Value = 0.1572327044025157
$ws.Range("R7").Value = 0.0880503144654088
$ws.Range("S7").Value = 0.3396226415094339

# Row 8
$ws.Range("B8").Value = 0.08385744234800839
$ws.Range("D8").Value = 0.01886792452830189
$ws.Range("F8").Value = 0.05870020964360587
$ws.Range("J8").Value = 0.0880503144654088
$ws.Range("O8").Value = 0.01257861635220126
$ws.Range("Q8").Value = 0.2222222222222222
$ws.Range("R8").Value = 0.05870020964360587
$ws.Range("S8").Value = 0.4570230607966457

# Row 9
$ws.Range("B9").Value = 0.1127450980392157
$ws.Range("D9").Value = 0.009803921568627451
$ws.Range("E9").Value = 0.004901960784313725
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("O9").Value = 0.0196078431372549
$ws.Range("Q9").Value = 0.1911764705882353
$ws.Range("R9").Value = 0.06862745098039216
$ws.Range("S9").Value = 0.4509803921568628

# Row 10
$ws.Range("B10").Value = 0.09365325077399381
$ws.Range("D10").Value = 0.02321981424148607
$ws.Range("F10").Value = 0.07507739938080496
$ws.Range("J10").Value = 0.119969040247678
$ws.Range("O10").Value = 0.01470588235294118
$ws.Range("Q10").Value = 0.2376160990712074
$ws.Range("R10").Value = 0.07894736842105263
$ws.Range("S10").Value = 0.3568111455108359

# Row 11
$ws.Range("G11").Value = 0.1558441558441558
$ws.Range("J11").Value = 0.1168831168831169
$ws.Range("K11").Value = 0.25
$ws.Range("L11").Value = 0.448051948051948
$ws.Range("S11").Value = 0.02922077922077922

# Row 12
$ws.Range("G12").Value = 0.673469387755102
$ws.Range("J12").Value = 0.2380952380952381
$ws.Range("K12").Value = 0.01360544217687075
$ws.Range("L12").Value = 0.04081632653061224
$ws.Range("S12").Value = 0.03401360544217687

# Row 15
$ws.Range("F15").Value = 0.02586206896551724
$ws.Range("H15").Value = 0.125
$ws.Range("I15").Value = 0.05603448275862069
$ws.Range("J15").Value = 0.3922413793103448
$ws.Range("K15").Value = 0.05603448275862069
$ws.Range("M15").Value = 0.008620689655172414
$ws.Range("O15").Value = 0.103448275862069
$ws.Range("S15").Value = 0.2327586206896552

# Row 16
$ws.Range("F16").Value = 0.03932584269662921
$ws.Range("H16").Value = 0.1853932584269663
$ws.Range("I16").Value = 0.101123595505618
$ws.Range("J16").Value = 0.3820224719101123
$ws.Range("K16").Value = 0.07303370786516854
$ws.Range("M16").Value = 0.01685393258426966
$ws.Range("O16").Value = 0.0449438202247191
$ws.Range("S16").Value = 0.1573033707865168

# Row 17
$ws.Range("F17").Value = 0.01764705882352941
$ws.Range("H17").Value = 0.1725490196078431
$ws.Range("I17").Value = 0.07843137254901961
$ws.Range("J17").Value = 0.4450980392156862
$ws.Range("K17").Value = 0.09803921568627451
$ws.Range("M17").Value = 0.007843137254901961
$ws.Range("O17").Value = 0.05882352941176471
$ws.Range("S17").Value = 0.1215686274509804

# Row 18
$ws.Range("F18").Value = 0.02339181286549707
$ws.Range("H18").Value = 0.2046783625730994
$ws.Range("I18").Value = 0.1169590643274854
$ws.Range("J18").Value = 0.3859649122807017
$ws.Range("K18").Value = 0.08771929824561403
$ws.Range("M18").Value = 0.005847953216374269
$ws.Range("O18").Value = 0.05847953216374269
$ws.Range("S18").Value = 0.1169590643274854

# Row 19
$ws.Range("F19").Value = 0.01755725190839695
$ws.Range("H19").Value = 0.2244274809160305
$ws.Range("I19").Value = 0.08778625954198473
$ws.Range("J19").Value = 0.3488549618320611
$ws.Range("K19").Value = 0.1015267175572519
$ws.Range("M19").Value = 0.0183206106870229
$ws.Range("O19").Value = 0.06870229007633588
$ws.Range("S19").Value = 0.132824427480916

Write-Host "Applied 101 cell updates across rows 2,3,6-12,15-19"
